# Update benchmark: 2025-11-01 06:35:04 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

# G2: was empty inline string, now "9 TL - 9 TL"
$ws.Range("G2").Value = "9 TL - 9 TL"

# G7: was empty inline string, now explanatory text
$ws.Range("G7").Value = "1 TRY (Kredi kartı ile ödemelerde ek olarak nakit avans faizi uygulanır.)"

# C13: update Azami (maximum) value
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"

# E13: update Azami (maximum) value
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 8.700 TL"
